## Latest Update with Validation - adds TC013-TC022 test case rows, updates
## TC011/TC012 wording, widens columns A/D/E, and moves the active
## selection/freeze-pane viewport down to the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column width adjustments (A, D, E get wider to fit the new, longer text)
# ColumnWidth (chars) maps to stored OOXML width via stored = ColumnWidth + 5/6,
# quantised to 1/6-character (6px) steps by this engine, so we pick the
# ColumnWidth value whose rounded result lands on the target stored width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 33.1666666666667   # -> stored width 34
$ws.Columns.Item(4).ColumnWidth = 55.3333333333333   # -> stored width ~56.1667 (target 56.140625)
$ws.Columns.Item(5).ColumnWidth = 43                 # -> stored width ~43.8333 (target 43.85546875)

# ---------------------------------------------------------------------------
# Row 32 (TC011_SavTypeMain_Add) - reword the "Add button disabled" result
# ---------------------------------------------------------------------------
$ws.Cells.Item(32, 1).Value = "TC011_SavTypeMain_Add"
$ws.Cells.Item(32, 2).Value = "Verify that user can add savings type with an empty data"
$ws.Cells.Item(32, 3).Value = "The admin has logged in on the system."
$ws.Cells.Item(32, 4).Value = "1. Click Add button."
$ws.Cells.Item(32, 6).Value = "1. Add button disabled."

# ---------------------------------------------------------------------------
# Row 33 (TC012_SavTypeMain_Add) - F33 turns from a stray numeric 1 into text,
# and a validation note is appended as new row 34
# ---------------------------------------------------------------------------
$ws.Cells.Item(33, 1).Value = "TC012_SavTypeMain_Add"
$ws.Cells.Item(33, 2).Value = "Verify that user can add numeric savings type name."
$ws.Cells.Item(33, 3).Value = "The admin has logged in on the system."
$ws.Cells.Item(33, 4).Value = "1. Enter the given savings type."
$ws.Cells.Item(33, 5).Value = "Savings Type Name: [Any numeric data]"
$ws.Cells.Item(33, 6).Value = "1. Savings Type Name textbox will show invalid icon."

$ws.Cells.Item(34, 4).Value = "2. Click Add button."
$ws.Cells.Item(34, 6).Value = "2. Add button disabled."

# ---------------------------------------------------------------------------
# New test case TC013_SavTypeMain_Add (rows 35-38)
# ---------------------------------------------------------------------------
$ws.Cells.Item(35, 1).Value = "TC013_SavTypeMain_Add"
$ws.Cells.Item(35, 2).Value = "Verify that user can add savings type with invalid data."
$ws.Cells.Item(35, 3).Value = "The admin has logged in on the system."
$ws.Cells.Item(35, 4).Value = "1. Enter the given data."
$ws.Cells.Item(35, 5).Value = "No. of Account Holders: [Alphanumeric data]"
$ws.Cells.Item(35, 6).Value = "1. Any given textboxes cannot input alphabetic characters except period."

$ws.Cells.Item(36, 5).Value = "Max. Withdrawal Amount: [Alphanumeric data]"
$ws.Cells.Item(37, 5).Value = "Maintaining Balance: [Alphanumeric data]"
$ws.Cells.Item(38, 5).Value = "Interest Rate: [Alphanumeric data]"

# ---------------------------------------------------------------------------
# New test case TC014_SavTypeMain_Update (rows 39-41)
# ---------------------------------------------------------------------------
$ws.Cells.Item(39, 1).Value = "TC014_SavTypeMain_Update"
$ws.Cells.Item(39, 2).Value = "Verify that user can change savings type."
$ws.Cells.Item(39, 3).Value = "The admin has logged in on the system."
$ws.Cells.Item(39, 4).Value = "1. Double click to select the savings type in the datagridview."
$ws.Cells.Item(39, 6).Value = "1.Dialog will appear for confirmation."

$ws.Cells.Item(40, 4).Value = "2. Change interest percentage to fixed"
$ws.Cells.Item(40, 6).Value = "2. The updated savings type will appear in the datagridview with update data."

$ws.Cells.Item(41, 4).Value = "3. Click Update button."

# ---------------------------------------------------------------------------
# New test case TC015_SavDormMain_Add (rows 42-44)
# ---------------------------------------------------------------------------
$ws.Cells.Item(42, 1).Value = "TC015_SavDormMain_Add"
$ws.Cells.Item(42, 2).Value = "Verify that user can add dormancy of a savings type"
$ws.Cells.Item(42, 3).Value = "The admin has logged in on the system."
$ws.Cells.Item(42, 4).Value = "1. Select Savings Type."

$ws.Cells.Item(43, 4).Value = "2. Input Inactive Duration and Amount Deducted"
$ws.Cells.Item(43, 5).Value = "Inactive Duration: 25"
$ws.Cells.Item(43, 6).Value = "1. Dialog box will appear. Record has been added."

$ws.Cells.Item(44, 4).Value = "3. Choose Fixed Amount and Active."
$ws.Cells.Item(44, 5).Value = "Amount Deducted: 100"
$ws.Cells.Item(44, 6).Value = "2. The added record will show in the datagridview."

# ---------------------------------------------------------------------------
# New test case TC016_SavDormMain_Update (rows 45-47)
# ---------------------------------------------------------------------------
$ws.Cells.Item(45, 1).Value = "TC016_SavDormMain_Update"
$ws.Cells.Item(45, 2).Value = "Verify that user can update dormancy of a savings type"
$ws.Cells.Item(45, 3).Value = "The admin has logged in on the system."
$ws.Cells.Item(45, 4).Value = "1. Double click the savings type that wants to be change."
$ws.Cells.Item(45, 5).Value = "Day to month"
$ws.Cells.Item(45, 6).Value = "1. Dialog box will appear. Record has been updated."

$ws.Cells.Item(46, 4).Value = "2. Change the inactivity period."

$ws.Cells.Item(47, 4).Value = "3. Click Update button."

# ---------------------------------------------------------------------------
# New test case TC017_TimeDepoTermRates_Add (rows 48-52)
# ---------------------------------------------------------------------------
$ws.Cells.Item(48, 1).Value = "TC017_TimeDepoTermRates_Add"
$ws.Cells.Item(48, 2).Value = "Verify that user can add terms and rates on the time deposit."
$ws.Cells.Item(48, 3).Value = "The admin has logged in on the system."
$ws.Cells.Item(48, 4).Value = "1. Enter valid data."
$ws.Cells.Item(48, 5).Value = "From: 1000.00"
$ws.Cells.Item(48, 6).Value = "1. Dialog box will appear. Record has been added."

$ws.Cells.Item(49, 4).Value = "2. Click Add button."
$ws.Cells.Item(49, 5).Value = "To: 4999.00"
$ws.Cells.Item(49, 6).Value = "2. The added record will show in the datagridview."

$ws.Cells.Item(50, 5).Value = "No. of Days: 60"
$ws.Cells.Item(51, 5).Value = "Interest Rate: 1.25"
$ws.Cells.Item(52, 5).Value = "Status: Active"

# ---------------------------------------------------------------------------
# New test case TC018_TimeDepoTermRates_Update (rows 53-55)
# ---------------------------------------------------------------------------
$ws.Cells.Item(53, 1).Value = "TC018_TimeDepoTermRates_Update"
$ws.Cells.Item(53, 2).Value = "Verify that user can update terms and rate on the time deposit."
$ws.Cells.Item(53, 3).Value = "The admin has logged in on the system."
$ws.Cells.Item(53, 4).Value = "1. Double click the terms and rates that you want to change."
$ws.Cells.Item(53, 5).Value = "Interest Rate: 0.25"
$ws.Cells.Item(53, 6).Value = "1. Dialog box will appear. Record has been updated."

$ws.Cells.Item(54, 4).Value = "2. Change Interest rate."
$ws.Cells.Item(54, 6).Value = "2. The updated terms and Rates will appear with updated data."

$ws.Cells.Item(55, 4).Value = "3. Click Update button."

# ---------------------------------------------------------------------------
# New test case TC019_TimeDepoPreTerm_Add (rows 56-57)
# ---------------------------------------------------------------------------
$ws.Cells.Item(56, 1).Value = "TC019_TimeDepoPreTerm_Add"
$ws.Cells.Item(56, 2).Value = "Verify that user can add pre termination penalty on the time deposit."
$ws.Cells.Item(56, 3).Value = "The admin has logged in on the system."
$ws.Cells.Item(56, 4).Value = "1. Enter valid data."

$ws.Cells.Item(57, 4).Value = "2. Click Add button."

# ---------------------------------------------------------------------------
# New test case TC020_MemberType_Add (rows 58-61)
# ---------------------------------------------------------------------------
$ws.Cells.Item(58, 1).Value = "TC020_MemberType_Add"
$ws.Cells.Item(58, 2).Value = "Verify that user can add Member type."
$ws.Cells.Item(58, 3).Value = "The admin has logged in on the system."
$ws.Cells.Item(58, 4).Value = "1. Enter a valid member type"
$ws.Cells.Item(58, 5).Value = "Member type: Another Member"
$ws.Cells.Item(58, 6).Value = "1. Dialog box will appear. Record has been added."

$ws.Cells.Item(59, 4).Value = "2. Click Add button."
$ws.Cells.Item(59, 5).Value = "Minimum Share: 100"
$ws.Cells.Item(59, 6).Value = "2. The added record will show in the datagridview."

$ws.Cells.Item(60, 5).Value = "Checked has certificate check box"
$ws.Cells.Item(61, 5).Value = "select active"

# ---------------------------------------------------------------------------
# New test case TC021_MemberType_Update (rows 62-64)
# ---------------------------------------------------------------------------
$ws.Cells.Item(62, 1).Value = "TC021_MemberType_Update"
$ws.Cells.Item(62, 2).Value = "Verify that user can update member type"
$ws.Cells.Item(62, 3).Value = "The admin has logged in on the system."
$ws.Cells.Item(62, 4).Value = "1. Double click the member type that you want to change."
$ws.Cells.Item(62, 5).Value = "Member type: Kiddie Member"
$ws.Cells.Item(62, 6).Value = "1. Dialog box will appear. Record had successfully updated."

$ws.Cells.Item(63, 4).Value = "2. Change Member type name."
$ws.Cells.Item(63, 6).Value = "2. the updated record will show in the datagridview."

$ws.Cells.Item(64, 4).Value = "3. Click Update button."

# ---------------------------------------------------------------------------
# New (still incomplete in the source workbook) test case TC022_ (row 65)
# ---------------------------------------------------------------------------
$ws.Cells.Item(65, 1).Value = "TC022_"

# ---------------------------------------------------------------------------
# View state: move the frozen-pane viewport and selection down to the new
# bottom of the sheet, and drop the D1 fixed top-left cell from the sheet view.
# ---------------------------------------------------------------------------
$ws.Range("A49").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 49
$ws.Range("A65").Select()
